$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 2385286.2
$ws.Range("I137").Value = 7144294
$ws.Range("K137").Value = 21432882
$ws.Range("M137").Value = -21430332

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Range("H52").Value = 40390
$ws.Range("J52").Value = 40390
$ws.Range("L52").Value = 40390
$ws.Range("N52").Value = -41026
# Row 61
$ws.Range("H61").Value = 43567864
$ws.Range("I61").Value = 66734892
$ws.Range("K61").Value = 66734892
$ws.Range("M61").Value = -66734680
# Row 74
$ws.Range("H74").Value = 9316412
$ws.Range("I74").Value = 12860106
$ws.Range("J74").Value = 102809.7
$ws.Range("K74").Value = 12860106
$ws.Range("L74").Value = 102809.7
$ws.Range("M74").Value = -12859232
$ws.Range("N74").Value = -104557.7
# Row 77
$ws.Range("H77").Value = 9316412
$ws.Range("I77").Value = 12860106
$ws.Range("J77").Value = 102809.7
$ws.Range("K77").Value = 64300530
$ws.Range("L77").Value = 514048.5
$ws.Range("M77").Value = -64296162
$ws.Range("N77").Value = -522784.5
# Row 122
$ws.Range("H122").Value = 9260882
$ws.Range("I122").Value = 1379.1111
$ws.Range("K122").Value = 4137.3333
$ws.Range("M122").Value = -1687.3333
# Row 132
$ws.Range("H132").Value = 76659.41
$ws.Range("I132").Value = 54837.21
$ws.Range("J132").Value = 128487.125
$ws.Range("K132").Value = 164511.63
$ws.Range("L132").Value = 385461.375
$ws.Range("M132").Value = -161981.63
$ws.Range("N132").Value = -390521.375
# Row 134
$ws.Range("H134").Value = 54507.688
$ws.Range("J134").Value = 54507.688
$ws.Range("L134").Value = 54507.688
$ws.Range("N134").Value = -64647.688
# Row 136
$ws.Range("H136").Value = 43567864
$ws.Range("I136").Value = 66734892
$ws.Range("K136").Value = 200204676
$ws.Range("M136").Value = -200202126

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 523.25
$ws.Range("I80").Value = 89
$ws.Range("J80").Value = 783.8
$ws.Range("K80").Value = 89
$ws.Range("L80").Value = 783.8
$ws.Range("M80").Value = 909
$ws.Range("N80").Value = -2779.8
# Row 83
$ws.Range("H83").Value = 523.25
$ws.Range("I83").Value = 89
$ws.Range("J83").Value = 783.8
$ws.Range("K83").Value = 445
$ws.Range("L83").Value = 3919
$ws.Range("M83").Value = 4547
$ws.Range("N83").Value = -13903
# Row 134
$ws.Range("H134").Value = 2002.7826
$ws.Range("I134").Value = 1922.5122
$ws.Range("J134").Value = 2661
$ws.Range("K134").Value = 5767.536599999999
$ws.Range("L134").Value = 7983
$ws.Range("M134").Value = -3232.536599999999
$ws.Range("N134").Value = -13053

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 46933.92
$ws.Range("I31").Value = 3031.3333
$ws.Range("J31").Value = 70176.47
$ws.Range("K31").Value = 3031.3333
$ws.Range("L31").Value = 70176.47
$ws.Range("M31").Value = -2736.3333
$ws.Range("N31").Value = -70766.47
# Row 34
$ws.Range("H34").Value = 46933.92
$ws.Range("I34").Value = 3031.3333
$ws.Range("J34").Value = 70176.47
$ws.Range("K34").Value = 3031.3333
$ws.Range("L34").Value = 70176.47
$ws.Range("M34").Value = -2829.3333
$ws.Range("N34").Value = -70580.47
# Row 58
$ws.Range("H58").Value = 34484100
$ws.Range("I58").Value = 45455932
$ws.Range("J58").Value = 1200
$ws.Range("K58").Value = 45455932
$ws.Range("L58").Value = 1200
$ws.Range("M58").Value = -45455729
$ws.Range("N58").Value = -1606
# Row 132
$ws.Range("H132").Value = 19171.736
$ws.Range("I132").Value = 1346.0212
$ws.Range("J132").Value = 102952.6
$ws.Range("K132").Value = 4038.063599999999
$ws.Range("L132").Value = 308857.8
$ws.Range("M132").Value = -1508.063599999999
$ws.Range("N132").Value = -313917.8
# Row 136
$ws.Range("H136").Value = 34484100
$ws.Range("I136").Value = 45455932
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 136367796
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -136365246
$ws.Range("N136").Value = -8700

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 2466.4
$ws.Range("I75").Value = 2110.6667
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 6332.000100000001
$ws.Range("L75").Value = 9000
$ws.Range("M75").Value = -5334.000100000001
$ws.Range("N75").Value = -10996
# Row 78
$ws.Range("H78").Value = 2466.4
$ws.Range("I78").Value = 2110.6667
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 18996.0003
$ws.Range("L78").Value = 27000
$ws.Range("M78").Value = -14004.0003
$ws.Range("N78").Value = -36984
# Row 125
$ws.Range("H125").Value = 2519
$ws.Range("J125").Value = 3611
$ws.Range("L125").Value = 10833
$ws.Range("N125").Value = -20673
# Row 131
$ws.Range("H131").Value = 952.95386
$ws.Range("J131").Value = 983.56665
$ws.Range("L131").Value = 2950.69995
$ws.Range("N131").Value = -13030.69995

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1177.7778
$ws.Range("J16").Value = 2900
$ws.Range("L16").Value = 2900
$ws.Range("N16").Value = -3240
# Row 93
$ws.Range("H93").Value = 1580.6
$ws.Range("I93").Value = 1350.75
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1350.75
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -102.75
$ws.Range("N93").Value = -4996
# Row 100
$ws.Range("H100").Value = 1673.1177
$ws.Range("I100").Value = 1405.8889
$ws.Range("K100").Value = 1405.8889
$ws.Range("M100").Value = -864.8888999999999
# Row 136
$ws.Range("H136").Value = 77485.63
$ws.Range("I136").Value = 53380.35
$ws.Range("J136").Value = 146357.86
$ws.Range("K136").Value = 160141.05
$ws.Range("L136").Value = 439073.58
$ws.Range("M136").Value = -157591.05
$ws.Range("N136").Value = -444173.58

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2667.7144
$ws.Range("I122").Value = 977.1667
$ws.Range("K122").Value = 2931.5001
$ws.Range("M122").Value = -481.5001000000002
# Row 132
$ws.Range("H132").Value = 33438.145
$ws.Range("I132").Value = 21663.959
$ws.Range("J132").Value = 73806.78999999999
$ws.Range("K132").Value = 64991.87699999999
$ws.Range("L132").Value = 221420.37
$ws.Range("M132").Value = -62461.87699999999
$ws.Range("N132").Value = -226480.37
# Row 136
$ws.Range("H136").Value = 46645.957
$ws.Range("I136").Value = 35532
$ws.Range("J136").Value = 66790
$ws.Range("K136").Value = 106596
$ws.Range("L136").Value = 200370
$ws.Range("M136").Value = -104046
$ws.Range("N136").Value = -205470
